$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (also refreshes the _FilterDatabase defined name automatically)
$ws.Name = "Produtos"

# Header row: "Categoria" -> "Tipo do Produto"
$ws.Range("D1").Value = "Tipo do Produto"

# Fix capitalization of the observation note
$ws.Range("G5").Value = "Conferir Estoque"

# Column G picks up column F's formatting/width (matches the saved file)
$ws.Range("F1:F25").Copy()
$ws.Range("G1:G25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the new selection left behind in the saved file
$ws.Range("G6").Select()
